$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Delete the "View Requirements.txt " bullet entirely.
# ------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -eq "View Requirements.txt `r") {
        $d.Paragraphs.Item($i).Range.Delete()
        break
    }
}

# ------------------------------------------------------------------
# 2) Insert the new ".env" bullet block right after the
#    "pip install <name of dependency>" line (before "Run app using
#    command:").
# ------------------------------------------------------------------
$anchorIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "*pip install*") {
        $anchorIndex = $i
        break
    }
}

$anchorPara = $d.Paragraphs.Item($anchorIndex)
$anchorPara.Range.InsertParagraphAfter() | Out-Null

# -- "Set up .env file for database use. " (top level bullet, ilvl 0)
$p = $d.Paragraphs.Item($anchorIndex + 1)
$p.Range.Text = "Set up .env file for database use. "
$p.Range.ListFormat.ListLevelNumber = 1

# -- "340DBHOST='classmysql.engr.oregonstate.edu'" (ilvl 1)
$p.Range.InsertParagraphAfter() | Out-Null
$p = $d.Paragraphs.Item($anchorIndex + 2)
$p.Range.Text = "340DBHOST='classmysql.engr.oregonstate.edu'"
$p.Range.ListFormat.ListLevelNumber = 2

# -- "340DBUSER= '<your db username here>'" (ilvl 1, 2 runs)
$p.Range.InsertParagraphAfter() | Out-Null
$p = $d.Paragraphs.Item($anchorIndex + 3)
$p.Range.Text = "340DBUSER= " + [char]0x2018 + "<your db username here>" + [char]0x2019
$p.Range.ListFormat.ListLevelNumber = 2

# -- "340DBPW='<your db password here>'" (ilvl 1, 3 runs)
$p.Range.InsertParagraphAfter() | Out-Null
$p = $d.Paragraphs.Item($anchorIndex + 4)
$p.Range.Text = "340DBPW=" + [char]0x2018 + "<your db password here>" + [char]0x2019
$p.Range.ListFormat.ListLevelNumber = 2

# -- "340DB='<your db name here>'" (ilvl 1, 3 runs)
$p.Range.InsertParagraphAfter() | Out-Null
$p = $d.Paragraphs.Item($anchorIndex + 5)
$p.Range.Text = "340DB=" + [char]0x2019 + "<your db name here>" + [char]0x2019
$p.Range.ListFormat.ListLevelNumber = 2

# ------------------------------------------------------------------
# 3) Swap the hard-coded port number for a placeholder.
# ------------------------------------------------------------------
$d.Content.Find.Execute("8525", $true, $false, $false, $false, $false,
                         $true, 1, $false, "<any port number>", 2) | Out-Null
